$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new dataset (ACS 2016 5-year DP05)
$ws.Name = "ACS_16_5YR_DP05_metadata"

# Full target content for columns A (variable code) and B (label) - every row here
# reflects the up-to-date ACS_16_5YR_DP05 metadata mapping.
$rows = @(
    @(1, "GEO.id", "Id"),
    @(2, "GEO.id2", "Id2"),
    @(3, "GEO.display-label", "Geography"),
    @(4, "HC01_VC03", "Estimate; SEX AND AGE - Total population"),
    @(8, "HC01_VC04", "Estimate; SEX AND AGE - Total population - Male"),
    @(12, "HC01_VC05", "Estimate; SEX AND AGE - Total population - Female"),
    @(16, "HC01_VC08", "Estimate; SEX AND AGE - Under 5 years"),
    @(20, "HC01_VC09", "Estimate; SEX AND AGE - 5 to 9 years"),
    @(24, "HC01_VC10", "Estimate; SEX AND AGE - 10 to 14 years"),
    @(28, "HC01_VC11", "Estimate; SEX AND AGE - 15 to 19 years"),
    @(32, "HC01_VC12", "Estimate; SEX AND AGE - 20 to 24 years"),
    @(36, "HC01_VC13", "Estimate; SEX AND AGE - 25 to 34 years"),
    @(40, "HC01_VC14", "Estimate; SEX AND AGE - 35 to 44 years"),
    @(44, "HC01_VC15", "Estimate; SEX AND AGE - 45 to 54 years"),
    @(48, "HC01_VC16", "Estimate; SEX AND AGE - 55 to 59 years"),
    @(52, "HC01_VC17", "Estimate; SEX AND AGE - 60 to 64 years"),
    @(56, "HC01_VC18", "Estimate; SEX AND AGE - 65 to 74 years"),
    @(60, "HC01_VC19", "Estimate; SEX AND AGE - 75 to 84 years"),
    @(64, "HC01_VC20", "Estimate; SEX AND AGE - 85 years and over"),
    @(68, "HC01_VC23", "Estimate; SEX AND AGE - Median age (years)"),
    @(72, "HC01_VC26", "Estimate; SEX AND AGE - 18 years and over"),
    @(76, "HC01_VC27", "Estimate; SEX AND AGE - 21 years and over"),
    @(80, "HC01_VC28", "Estimate; SEX AND AGE - 62 years and over"),
    @(84, "HC01_VC29", "Estimate; SEX AND AGE - 65 years and over"),
    @(88, "HC01_VC32", "Estimate; SEX AND AGE - 18 years and over"),
    @(92, "HC01_VC33", "Estimate; SEX AND AGE - 18 years and over - Male"),
    @(96, "HC01_VC34", "Estimate; SEX AND AGE - 18 years and over - Female"),
    @(100, "HC01_VC37", "Estimate; SEX AND AGE - 65 years and over"),
    @(104, "HC01_VC38", "Estimate; SEX AND AGE - 65 years and over - Male"),
    @(108, "HC01_VC39", "Estimate; SEX AND AGE - 65 years and over - Female"),
    @(112, "HC01_VC43", "Estimate; RACE - Total population"),
    @(116, "HC01_VC44", "Estimate; RACE - Total population - One race"),
    @(120, "HC01_VC45", "Estimate; RACE - Total population - Two or more races"),
    @(124, "HC01_VC48", "Estimate; RACE - One race"),
    @(128, "HC01_VC49", "Estimate; RACE - One race - White"),
    @(132, "HC01_VC50", "Estimate; RACE - One race - Black or African American"),
    @(136, "HC01_VC51", "Estimate; RACE - One race - American Indian and Alaska Native"),
    @(140, "HC01_VC52", "Estimate; RACE - One race - American Indian and Alaska Native - Cherokee tribal grouping"),
    @(144, "HC01_VC53", "Estimate; RACE - One race - American Indian and Alaska Native - Chippewa tribal grouping"),
    @(148, "HC01_VC54", "Estimate; RACE - One race - American Indian and Alaska Native - Navajo tribal grouping"),
    @(152, "HC01_VC55", "Estimate; RACE - One race - American Indian and Alaska Native - Sioux tribal grouping"),
    @(156, "HC01_VC56", "Estimate; RACE - One race - Asian"),
    @(160, "HC01_VC57", "Estimate; RACE - One race - Asian - Asian Indian"),
    @(164, "HC01_VC58", "Estimate; RACE - One race - Asian - Chinese"),
    @(168, "HC01_VC59", "Estimate; RACE - One race - Asian - Filipino"),
    @(172, "HC01_VC60", "Estimate; RACE - One race - Asian - Japanese"),
    @(176, "HC01_VC61", "Estimate; RACE - One race - Asian - Korean"),
    @(180, "HC01_VC62", "Estimate; RACE - One race - Asian - Vietnamese"),
    @(184, "HC01_VC63", "Estimate; RACE - One race - Asian - Other Asian"),
    @(188, "HC01_VC64", "Estimate; RACE - One race - Native Hawaiian and Other Pacific Islander"),
    @(192, "HC01_VC65", "Estimate; RACE - One race - Native Hawaiian and Other Pacific Islander - Native Hawaiian"),
    @(196, "HC01_VC66", "Estimate; RACE - One race - Native Hawaiian and Other Pacific Islander - Guamanian or Chamorro"),
    @(200, "HC01_VC67", "Estimate; RACE - One race - Native Hawaiian and Other Pacific Islander - Samoan"),
    @(204, "HC01_VC68", "Estimate; RACE - One race - Native Hawaiian and Other Pacific Islander - Other Pacific Islander"),
    @(208, "HC01_VC69", "Estimate; RACE - One race - Some other race"),
    @(212, "HC01_VC70", "Estimate; RACE - Two or more races"),
    @(216, "HC01_VC71", "Estimate; RACE - Two or more races - White and Black or African American"),
    @(220, "HC01_VC72", "Estimate; RACE - Two or more races - White and American Indian and Alaska Native"),
    @(224, "HC01_VC73", "Estimate; RACE - Two or more races - White and Asian"),
    @(228, "HC01_VC74", "Estimate; RACE - Two or more races - Black or African American and American Indian and Alaska Native"),
    @(232, "HC01_VC77", "Estimate; RACE - Race alone or in combination with one or more other races - Total population"),
    @(236, "HC01_VC78", "Estimate; RACE - Race alone or in combination with one or more other races - Total population - White"),
    @(240, "HC01_VC79", "Estimate; RACE - Race alone or in combination with one or more other races - Total population - Black or African American"),
    @(244, "HC01_VC80", "Estimate; RACE - Race alone or in combination with one or more other races - Total population - American Indian and Alaska Native"),
    @(248, "HC01_VC81", "Estimate; RACE - Race alone or in combination with one or more other races - Total population - Asian"),
    @(252, "HC01_VC82", "Estimate; RACE - Race alone or in combination with one or more other races - Total population - Native Hawaiian and Other Pacific Islander"),
    @(256, "HC01_VC83", "Estimate; RACE - Race alone or in combination with one or more other races - Total population - Some other race"),
    @(260, "HC01_VC87", "Estimate; HISPANIC OR LATINO AND RACE - Total population"),
    @(264, "HC01_VC88", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Hispanic or Latino (of any race)"),
    @(268, "HC01_VC89", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Hispanic or Latino (of any race) - Mexican"),
    @(272, "HC01_VC90", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Hispanic or Latino (of any race) - Puerto Rican"),
    @(276, "HC01_VC91", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Hispanic or Latino (of any race) - Cuban"),
    @(280, "HC01_VC92", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Hispanic or Latino (of any race) - Other Hispanic or Latino"),
    @(284, "HC01_VC93", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino"),
    @(288, "HC01_VC94", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - White alone"),
    @(292, "HC01_VC95", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Black or African American alone"),
    @(296, "HC01_VC96", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - American Indian and Alaska Native alone"),
    @(300, "HC01_VC97", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Asian alone"),
    @(304, "HC01_VC98", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Native Hawaiian and Other Pacific Islander alone"),
    @(308, "HC01_VC99", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Some other race alone"),
    @(312, "HC01_VC100", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Two or more races"),
    @(316, "HC01_VC101", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Two or more races - Two races including Some other race"),
    @(320, "HC01_VC102", "Estimate; HISPANIC OR LATINO AND RACE - Total population - Not Hispanic or Latino - Two or more races - Two races excluding Some other race, and Three or more races"),
    @(324, "HC01_VC104", "Estimate; HISPANIC OR LATINO AND RACE - Total housing units"),
    @(328, "HC01_VC108", "Estimate; CITIZEN, VOTING AGE POPULATION - Citizen, 18 and over population"),
    @(332, "HC01_VC109", "Estimate; CITIZEN, VOTING AGE POPULATION - Citizen, 18 and over population - Male"),
    @(336, "HC01_VC110", "Estimate; CITIZEN, VOTING AGE POPULATION - Citizen, 18 and over population - Female")
)

foreach ($row in $rows) {
    $r = $row[0]
    $colA = $row[1]
    $colB = $row[2]
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
}

# Column widths to fit the new (wider) label text
$ws.Columns.Item(1).ColumnWidth = 15.6640625
$ws.Columns.Item(2).ColumnWidth = 157.1640625

# Restore a sensible selection near the bottom of the refreshed data
$ws.Range("B339").Select() | Out-Null
